$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = "thêm API update avatar"
$ws.Range("C9").Value = "user"
$ws.Range("D9").Value = "trung bình"
$ws.Range("E9").Value = "chờ"

$ws.Range("B10").Select()
